# 9.5.2 indicator sheet update:
#  - Replace the outdated Russian header in B1 with the new translation
#    ("Численность специалистов-исследователей..." -> "Количество
#    исследователей (в эквиваленте полной занятости)...").
#  - Add a new year column (Q) with 2023 data: header 2023 and the value
#    631, matching the formatting already used by the neighbouring 2022
#    column (P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the Russian header text in B1.
$ws.Range("B1").Value = "9.5.2 Количество исследователей (в эквиваленте полной занятости) на миллион жителей"

# 2) Add the 2023 column, copying formatting from the 2022 column (P)
#    so the new cells share the same number format/font/borders.
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 2023

$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 631
